# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

$ws.Range("A2").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B2").Value = "C01号直流"
$ws.Range("C2").Value = "2025-01-25 13:46:36"
$ws.Range("D2").Value = 46002.28700231481

$ws.Range("A3").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B3").Value = "C02号直流"
$ws.Range("C3").Value = "2025-01-25 17:13:47"
$ws.Range("D3").Value = 46002.28700231481

$ws.Range("A4").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B4").Value = "C03号直流"
$ws.Range("C4").Value = "2025-01-25 14:14:24"
$ws.Range("D4").Value = 46002.28700231481

$ws.Range("A5").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B5").Value = "C04号直流"
$ws.Range("C5").Value = "2025-01-25 06:24:40"
$ws.Range("D5").Value = 46002.28700231481

$ws.Range("A6").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B6").Value = "C05号直流"
$ws.Range("C6").Value = "2025-01-25 16:01:40"
$ws.Range("D6").Value = 46002.28700231481

$ws.Range("A7").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B7").Value = "D01号直流"
$ws.Range("C7").Value = "2025-01-25 18:30:24"
$ws.Range("D7").Value = 46002.28700231481

$ws.Range("A8").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B8").Value = "D02号直流"
$ws.Range("C8").Value = "2025-01-25 15:39:19"
$ws.Range("D8").Value = 46002.28700231481

$ws.Range("A9").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B9").Value = "D03号直流"
$ws.Range("C9").Value = "2025-01-25 16:09:35"
$ws.Range("D9").Value = 46002.28700231481

$ws.Range("A10").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B10").Value = "D04号直流"
$ws.Range("C10").Value = "2025-01-25 18:29:02"
$ws.Range("D10").Value = 46002.28700231481

$ws.Range("A11").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B11").Value = "D05号直流"
$ws.Range("C11").Value = "2025-01-25 18:27:29"
$ws.Range("D11").Value = 46002.28700231481

$ws.Range("A12").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B12").Value = "E01号直流"
$ws.Range("C12").Value = "2025-01-25 15:22:58"
$ws.Range("D12").Value = 46002.28700231481

$ws.Range("A13").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B13").Value = "E02号直流"
$ws.Range("C13").Value = "2025-01-25 16:45:57"
$ws.Range("D13").Value = 46002.28700231481

$ws.Range("A14").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B14").Value = "E03号直流"
$ws.Range("C14").Value = "2025-01-25 02:54:59"
$ws.Range("D14").Value = 46002.28700231481

$ws.Range("A15").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B15").Value = "E04号直流"
$ws.Range("C15").Value = "2025-01-25 17:08:37"
$ws.Range("D15").Value = 46002.28700231481

$ws.Range("A16").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B16").Value = "004B号直流"
$ws.Range("C16").Value = "2025-02-19 00:26:27"
$ws.Range("D16").Value = 46002.28700231481

$ws.Range("A17").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B17").Value = "701号直流"
$ws.Range("C17").Value = 45927.457337962966
$ws.Range("D17").Value = 46002.28700231481

$ws.Range("A18").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B18").Value = "007B号直流"
$ws.Range("C18").Value = "2025-08-10 13:17:02"
$ws.Range("D18").Value = 46002.28700231481

$ws.Range("A19").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B19").Value = "101号直流"
$ws.Range("C19").Value = 45987.55260416667
$ws.Range("D19").Value = 46002.28700231481

$ws.Range("A20").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B20").Value = "702号直流"
$ws.Range("C20").Value = 45997.07376157407
$ws.Range("D20").Value = 46002.28700231481

$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "603号直流"
$ws.Range("C21").Value = 45999.582453703704
$ws.Range("D21").Value = 46002.28700231481

$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "705号直流"
$ws.Range("C22").Value = 46000.51550925926
$ws.Range("D22").Value = 46002.28700231481

$ws.Range("A23").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B23").Value = "703号直流"
$ws.Range("C23").Value = 46000.52245370371
$ws.Range("D23").Value = 46002.28700231481

$ws.Range("A24").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B24").Value = "211号直流"
$ws.Range("C24").Value = 46000.55601851852
$ws.Range("D24").Value = 46002.28700231481

$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "B01号直流"
$ws.Range("C25").Value = 46000.56130787037
$ws.Range("D25").Value = 46002.28700231481

$ws.Range("A26").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B26").Value = "103号直流"
$ws.Range("C26").Value = 46000.57625
$ws.Range("D26").Value = 46002.28700231481

$ws.Range("A27").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B27").Value = "401号直流"
$ws.Range("C27").Value = 46001.05726851852
$ws.Range("D27").Value = 46002.28700231481

$ws.Range("A28").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B28").Value = "404号直流"
$ws.Range("C28").Value = 46001.08253472222
$ws.Range("D28").Value = 46002.28700231481

$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "504号直流"
$ws.Range("C29").Value = 46001.228842592594
$ws.Range("D29").Value = 46002.28700231481

$ws.Range("A30").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B30").Value = "201号直流"
$ws.Range("C30").Value = 46001.242847222224
$ws.Range("D30").Value = 46002.28700231481

$ws.Range("A31").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B31").Value = "110号直流"
$ws.Range("C31").Value = 46001.43746527778
$ws.Range("D31").Value = 46002.28700231481

$ws.Range("A32").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B32").Value = "003B号直流"
$ws.Range("C32").Value = 46001.51935185185
$ws.Range("D32").Value = 46002.28700231481

$ws.Range("A33").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B33").Value = "101号直流"
$ws.Range("C33").Value = 46001.53508101852
$ws.Range("D33").Value = 46002.28700231481

$ws.Range("A34").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B34").Value = "102号直流"
$ws.Range("C34").Value = 46001.542766203704
$ws.Range("D34").Value = 46002.28700231481

$ws.Range("A35").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B35").Value = "805号直流"
$ws.Range("C35").Value = 46001.55304398148
$ws.Range("D35").Value = 46002.28700231481

$ws.Range("A36").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B36").Value = "905号直流"
$ws.Range("C36").Value = 46001.55494212963
$ws.Range("D36").Value = 46002.28700231481

$ws.Range("A37").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B37").Value = "104号直流"
$ws.Range("C37").Value = 46001.55981481481
$ws.Range("D37").Value = 46002.28700231481

$ws.Range("A38").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B38").Value = "006A号直流"
$ws.Range("C38").Value = 46001.56123842593
$ws.Range("D38").Value = 46002.28700231481

$ws.Range("A39").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B39").Value = "108号直流"
$ws.Range("C39").Value = 46001.56350694445
$ws.Range("D39").Value = 46002.28700231481

$ws.Range("A40").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B40").Value = "002B号直流"
$ws.Range("C40").Value = 46001.5640625
$ws.Range("D40").Value = 46002.28700231481

$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "401号直流"
$ws.Range("C41").Value = 46001.578101851854
$ws.Range("D41").Value = 46002.28700231481

$ws.Range("A42").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B42").Value = "105号直流"
$ws.Range("C42").Value = 46001.58188657407
$ws.Range("D42").Value = 46002.28700231481

$ws.Range("A43").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B43").Value = "008B号直流"
$ws.Range("C43").Value = 46001.58571759259
$ws.Range("D43").Value = 46002.28700231481

$ws.Range("A44").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B44").Value = "305号直流"
$ws.Range("C44").Value = 46001.591412037036
$ws.Range("D44").Value = 46002.28700231481

$ws.Range("A45").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B45").Value = "209号直流"
$ws.Range("C45").Value = 46001.59248842593
$ws.Range("D45").Value = 46002.28700231481

$ws.Range("A46").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B46").Value = "902号直流"
$ws.Range("C46").Value = 46001.615648148145
$ws.Range("D46").Value = 46002.28700231481

$ws.Range("A47").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B47").Value = "A01号直流"
$ws.Range("C47").Value = 46001.61832175926
$ws.Range("D47").Value = 46002.28700231481

$ws.Range("A48").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B48").Value = "107号直流"
$ws.Range("C48").Value = 46001.63171296296
$ws.Range("D48").Value = 46002.28700231481

$ws.Range("A49").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B49").Value = "904号直流"
$ws.Range("C49").Value = 46001.65274305556
$ws.Range("D49").Value = 46002.28700231481

$ws.Range("A50").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B50").Value = "103号直流"
$ws.Range("C50").Value = 46001.68547453704
$ws.Range("D50").Value = 46002.28700231481

$ws.Range("A51").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B51").Value = "402号直流"
$ws.Range("C51").Value = 46001.71864583333
$ws.Range("D51").Value = 46002.28700231481

$ws.Range("A52").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B52").Value = "207号直流"
$ws.Range("C52").Value = 46001.74553240741
$ws.Range("D52").Value = 46002.28700231481

$ws.Range("A53").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B53").Value = "111号直流"
$ws.Range("C53").Value = 46001.74832175926
$ws.Range("D53").Value = 46002.28700231481

$ws.Range("A54").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B54").Value = "603号直流"
$ws.Range("C54").Value = 46001.76337962963
$ws.Range("D54").Value = 46002.28700231481

$ws.Range("H8").Select()
